{"js": "// Update the date title and the 25 division problems in the table.\n// Replacements are applied strictly by position (paragraph order / row-col\n// index) because several of the old values (e.g. \"37\u00f78=\") repeat with\n// different replacements, so a naive global text search-and-replace would\n// be ambiguous.\n\nconst titleReplacement = { from: \"2025-09-29 Monday\", to: \"2025-09-30 Tuesday\" };\n\n// Row-major order of the 25 visible problem cells (5 rows x 5 cols).\nconst cellReplacements = [\n  [\"27\u00f76=\", \"71\u00f79=\"], [\"41\u00f72=\", \"50\u00f76=\"], [\"51\u00f79=\", \"81\u00f78=\"], [\"37\u00f75=\", \"46\u00f72=\"], [\"66\u00f72=\", \"21\u00f78=\"],\n  [\"29\u00f72=\", \"46\u00f75=\"], [\"43\u00f76=\", \"95\u00f78=\"], [\"88\u00f74=\", \"82\u00f74=\"], [\"82\u00f76=\", \"54\u00f79=\"], [\"70\u00f75=\", \"90\u00f72=\"],\n  [\"57\u00f77=\", \"25\u00f75=\"], [\"16\u00f76=\", \"97\u00f78=\"], [\"37\u00f78=\", \"83\u00f78=\"], [\"58\u00f78=\", \"81\u00f74=\"], [\"40\u00f73=\", \"49\u00f78=\"],\n  [\"38\u00f77=\", \"87\u00f78=\"], [\"80\u00f74=\", \"73\u00f73=\"], [\"12\u00f73=\", \"11\u00f72=\"], [\"37\u00f78=\", \"57\u00f73=\"], [\"94\u00f72=\", \"90\u00f77=\"],\n  [\"63\u00f79=\", \"86\u00f76=\"], [\"65\u00f73=\", \"52\u00f73=\"], [\"85\u00f74=\", \"68\u00f72=\"], [\"68\u00f75=\", \"93\u00f75=\"], [\"58\u00f73=\", \"27\u00f78=\"],\n];\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph -------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nif (titlePara.text.trim() === titleReplacement.from) {\n  titlePara.insertText(titleReplacement.to, \"Replace\");\n}\n\n// --- 2. Table cells -------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values; // 2D array [row][col] of cell text\nlet k = 0; // index into cellReplacements, filled row-major, skipping blanks\nfor (let r = 0; r < values.length; r++) {\n  for (let c = 0; c < values[r].length; c++) {\n    const cellText = (values[r][c] || \"\").trim();\n    if (cellText === \"\") continue;\n    if (k >= cellReplacements.length) continue;\n    const [from, to] = cellReplacements[k];\n    k++;\n    if (cellText === from) {\n      table.getCell(r, c).value = to;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 division problems in the table.\n# Replacements are applied strictly by position (paragraph / table row-col\n# index) because several of the old values (e.g. \"37\u00f78=\") repeat with\n# different replacements, so a naive global Find/Replace would be ambiguous.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph -------------------------------------------------\n$titleFrom = \"2025-09-29 Monday\"\n$titleTo = \"2025-09-30 Tuesday\"\n\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd(\"`r\") -eq $titleFrom) {\n    $titlePara.Range.Text = $titleTo\n}\n\n# --- 2. Table cells -------------------------------------------------\n# Row-major order of the 25 visible problem cells (5 rows x 5 cols).\n$cellFrom = @(\n    \"27\u00f76=\", \"41\u00f72=\", \"51\u00f79=\", \"37\u00f75=\", \"66\u00f72=\",\n    \"29\u00f72=\", \"43\u00f76=\", \"88\u00f74=\", \"82\u00f76=\", \"70\u00f75=\",\n    \"57\u00f77=\", \"16\u00f76=\", \"37\u00f78=\", \"58\u00f78=\", \"40\u00f73=\",\n    \"38\u00f77=\", \"80\u00f74=\", \"12\u00f73=\", \"37\u00f78=\", \"94\u00f72=\",\n    \"63\u00f79=\", \"65\u00f73=\", \"85\u00f74=\", \"68\u00f75=\", \"58\u00f73=\"\n)\n$cellTo = @(\n    \"71\u00f79=\", \"50\u00f76=\", \"81\u00f78=\", \"46\u00f72=\", \"21\u00f78=\",\n    \"46\u00f75=\", \"95\u00f78=\", \"82\u00f74=\", \"54\u00f79=\", \"90\u00f72=\",\n    \"25\u00f75=\", \"97\u00f78=\", \"83\u00f78=\", \"81\u00f74=\", \"49\u00f78=\",\n    \"87\u00f78=\", \"73\u00f73=\", \"11\u00f72=\", \"57\u00f73=\", \"90\u00f77=\",\n    \"86\u00f76=\", \"52\u00f73=\", \"68\u00f72=\", \"93\u00f75=\", \"27\u00f78=\"\n)\n\n$tbl = $d.Tables.Item(1)\n$k = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd(\"`a\").TrimEnd(\"`r\")\n        if ($cellText -eq \"\") { continue }\n        if ($k -ge $cellFrom.Length) { continue }\n        if ($cellText -eq $cellFrom[$k]) {\n            $cell.Range.Text = $cellTo[$k]\n        }\n        $k++\n    }\n}\n"}
